# "added and reworked MC questions"
# Remove the two rows for desc_stat-10.Rnw and desc_stat-11.Rnw (rows 11-12),
# which shifts every following row up by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two question rows (file_name / points) for desc_stat-10.Rnw and
# desc_stat-11.Rnw. They sit at rows 11 and 12 right now.
$ws.Range("A11:A12").EntireRow.Delete()

# The conditional formatting on the points column still targets the old
# (now too-large) range; re-point it at the shrunk data range.
$fcs = $ws.Range("B2:B83").FormatConditions
$fcCount = $fcs.Count()
for ($i = 1; $i -le $fcCount; $i++) {
    $fc = $fcs.Item($i)
    $fc.ModifyAppliesToRange($ws.Range("B2:B81"))
}

# Restore the view: scroll back to the top and select E18 (matches the
# post-edit author session, rather than the stale E40/topLeftCell=A24).
[void]$ws.Range("A1").Select()
[void]$ws.Range("E18").Select()
